$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shop data excel fix: each upgrade-shop row now references its own
# unique item id instead of all rows sharing 5001.
$ws.Range("B5").Value = 5002
$ws.Range("B6").Value = 5003
$ws.Range("B7").Value = 5004
$ws.Range("B8").Value = 5005
$ws.Range("B9").Value = 5006
$ws.Range("B10").Value = 5007

$ws.Range("E7").Select()
